$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two applicant rows that previously held "Dhanush " / "bro"
$ws.Range("A4").Value = "aa"
$ws.Range("A5").Value = "asmi"

# Add a new row 6: an (empty-but-hyperlink-styled) B6 cell and a score in C6
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("C6").Value = 77

# Widen column B (Email column) to fit - matches saved width of 18.5
$ws.Columns.Item(2).ColumnWidth = 17 + 2/3

# Move/leave the active selection on A6, like in the edited file
$ws.Range("A6").Select()
